# map for mycat to LVIS_cat
# Append the LVIS category list (rows 4-26, column A) below the existing
# My Category / LVIS Category mapping header + two sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Identity",
    "Home interior",
    "Vehicle Plate",
    "Bystander",
    "Food",
    "Paper/Document/Label",
    "Screen",
    "Clothing",
    "Scenery",
    "Pet",
    "Book",
    "Photo",
    "Machine",
    "Table",
    "Electronic devices",
    "Toiletries",
    "Toy",
    "Finger",
    "Cigarettes",
    "Accident",
    "Music instrument",
    "Nudity",
    "Accessory"
)

$startRow = 4
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}

# "Finger" (row 21) carries a distinct font (same face/size, just applied
# directly rather than through the sheet default).
$ws.Range("A21").Font.Name = "游ゴシック"

[void]$ws.Range("B11").Select()
